$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.137.30"
$ws.Range("E2").Value = "  +0.01%  "

# Row 3
$ws.Range("D3").Value = "3.125.40"
$ws.Range("E3").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").Value = "579.59"
$ws.Range("E5").Value = "  -0.04%  "

# Row 6
$ws.Range("D6").Value = "177.97"
$ws.Range("E6").Value = "  +2.19%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").Value = "3.124.22"
$ws.Range("E8").Value = "  +0.31%  "

# Row 9
$ws.Range("E9").Value = "  -1.21%  "

# Row 10
$ws.Range("D10").Value = "6.42"
$ws.Range("E10").Value = "  -0.52%  "

# Row 11
$ws.Range("E11").Value = "  -1.58%  "

# Row 12
$ws.Range("E12").Value = "  -1.07%  "

# Row 13
$ws.Range("E13").Value = "  -2.34%  "

# Row 14
$ws.Range("D14").Value = "36.46"
$ws.Range("E14").Value = "  -1.62%  "

# Row 15
$ws.Range("E15").Value = "  -0.34%  "

# Row 16
$ws.Range("D16").Value = "3.648.24"
$ws.Range("E16").Value = "  +0.33%  "

# Row 17
$ws.Range("D17").Value = "67.054.56"
$ws.Range("E17").Value = "  -0.09%  "

# Row 18
$ws.Range("D18").Value = "7.05"
$ws.Range("E18").Value = "  -0.93%  "

# Row 19
$ws.Range("D19").Value = "17.01"
$ws.Range("E19").Value = "  +2.03%  "

# Row 20
$ws.Range("D20").Value = "3.127.36"
$ws.Range("E20").Value = "  +0.23%  "

# Row 21
$ws.Range("D21").Value = "489.35"
$ws.Range("E21").Value = "  -0.58%  "

# Row 22
$ws.Range("D22").Value = "7.81"
$ws.Range("E22").Value = "  -1.32%  "

# Row 23
$ws.Range("D23").Value = "0.697"
$ws.Range("E23").Value = "  -1.38%  "

# Row 24
$ws.Range("E24").Value = "  -0.39%  "

# Row 25
$ws.Range("E25").Value = "  -3.03%  "

# Row 26
$ws.Range("E26").Value = "  -0.98%  "

# Row 27
$ws.Range("D27").Value = "10.34"
$ws.Range("E27").Value = "  -1.29%  "

# Row 28
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.01%  "

# Row 29
$ws.Range("E29").Value = "  +1.75%  "

# Row 30
$ws.Range("E30").Value = "  -1.92%  "

# Row 31
$ws.Range("E31").Value = "  -2.23%  "

# Row 32
$ws.Range("D32").Value = "28.24"
$ws.Range("E32").Value = "  -1.10%  "

# Row 33
$ws.Range("E33").Value = "  -0.86%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0948"
$ws.Range("E34").Value = "  +0.06%  "

# Row 35
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("D36").Value = "48.84"
$ws.Range("E36").Value = "  +2.98%  "

# Row 37
$ws.Range("D37").Value = "5.67"
$ws.Range("E37").Value = "  -3.64%  "

# Row 38
$ws.Range("D38").Value = "0.948"
$ws.Range("E38").Value = "  -2.82%  "

# Row 39
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "0.313"
$ws.Range("E39").Value = "  +0.69%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "49.45"
$ws.Range("E40").Value = "  -1.32%  "

# Row 41
$ws.Range("D41").Value = "2.01"
$ws.Range("E41").Value = "  -2.11%  "

# Row 42
$ws.Range("E42").Value = "  -0.20%  "

# Row 43
$ws.Range("E43").Value = "  -1.48%  "

# Row 44
$ws.Range("E44").Value = "  +3.35%  "

# Row 45
$ws.Range("D45").Value = "2.799.17"
$ws.Range("E45").Value = "  -0.82%  "

# Row 46
$ws.Range("D46").Value = "376.14"

# Row 47
$ws.Range("D47").Value = "0.0349"
$ws.Range("E47").Value = "  -1.16%  "

# Row 48
$ws.Range("D48").Value = "135.12"
$ws.Range("E48").Value = "  -0.35%  "

# Row 49
$ws.Range("E49").Value = "  -0.03%  "

# Row 50
$ws.Range("D50").Value = "25.15"
$ws.Range("E50").Value = "  +1.10%  "

# Row 51
$ws.Range("E51").Value = "  +2.28%  "
